# The <id>p030r_2</id> tag in the document is split across three runs:
#   1. "<id>"   - Courier New, color 7f6000, sz 18
#   2. "p030r_2" - default run formatting
#   3. "</id>"  - Courier New, color 7f6000, sz 18
#
# The edit collapses these into a single run containing the full text
# "<id>p030r_2</id>" (keeping the Courier New / 7f6000 / sz18 look of the
# surrounding tag runs). Doing a Find/Replace over the exact same text
# (which is unique in the document) makes Word re-emit the matched span
# as one run, which merges the three runs into one as required.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p030r_2</id>",  # FindText
    $true,               # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "<id>p030r_2</id>",   # ReplaceWith
    2                      # Replace (wdReplaceAll)
)
